$wb = $excel.ActiveWorkbook

# --- "Initial values" sheet ---
$ws1 = $wb.Worksheets.Item("Initial values")

# B2: 350/0.9 -> 350*1000/0.9
$ws1.Range("B2").Formula = "=350*1000/0.9"

# B8: ROUND(0.03/0.9,4) -> ROUND(0.03*1000/0.9,4)
$ws1.Range("B8").Formula = "=ROUND(0.03*1000/0.9,4)"

# B10: 1 -> 1000
$ws1.Range("B10").Value = 1000

# Update the active selection on this sheet from D3 to B3
$ws1.Activate()
$ws1.Range("B3").Select()

# --- "Advancements1" sheet ---
$ws4 = $wb.Worksheets.Item("Advancements1")

# B5 loses its style (numFmtId=0 applyNumberFormat xf) -> default style
$ws4.Range("B5").Style = "Normal"
